# Apply the data refresh described in the commit:
#  - TestID (D) and TestName (E) updated for every data row (2-22)
#  - Total Marks (R) reduced by the old Essay 4 (AF) score for rows that had one
#  - Rank (T) recomputed for the affected rows after the Total Marks change
#  - Essay 4 (AF) score cleared for rows that had one

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TestID / TestName change applies uniformly to every data row (2..22) ---
$newTestId = "7a6eafab-5fa3-4f6b-9220-96f64351c8b8"
$newTestName = "IBA Mock 1"

$ws.Range("D2:D22").Value = $newTestId
$ws.Range("E2:E22").Value = $newTestName

# --- Per-row Total Marks / Rank / Essay 4 updates ---
# row -> (new Total Marks, new Rank or $null if unchanged)
$updates = @{
    2  = @{ R = 22.5;  T = $null }
    4  = @{ R = 44.5;  T = $null }
    6  = @{ R = 38.25; T = 9 }
    7  = @{ R = 41.25; T = 7 }
    9  = @{ R = 38.75; T = 8 }
    11 = @{ R = 50.25; T = 1 }
    15 = @{ R = 47.75; T = 4 }
    16 = @{ R = 48.25; T = 2 }
    18 = @{ R = 29.25; T = $null }
    19 = @{ R = 48;    T = $null }
    20 = @{ R = 34;    T = $null }
    22 = @{ R = 42.25; T = 6 }
}

foreach ($row in $updates.Keys) {
    $u = $updates[$row]

    $ws.Range("R$row").Value = $u.R
    if ($null -ne $u.T) {
        $ws.Range("T$row").Value = $u.T
    }

    # Essay 4 score is removed now that it has been folded into Total Marks
    $ws.Range("AF$row").ClearContents()
}
